# Change RH/temp sensor to STH21 (HTU21 -> SHT21) on the "command" sheet,
# row 22: update part name, Farnell ref + unit price, which ripples through
# the existing SUM/shared formulas, and refresh the sheet's selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("command")
$ws.Activate()

# Row 22 - swap the sensor reference and fill in its Farnell part number /
# unit price (previously blank placeholders).
$ws.Range("A22").Value = "SHT21"
$ws.Range("D22").Value = 1855468
$ws.Range("E22").Value = 3.79

# Move the active selection to reflect where the edit was made.
$ws.Range("E22").Select()
